$d = $word.ActiveDocument

# 1. Add new text after "chỉ 1 cổng active tại 1 thời điểm."
$d.Content.Find.Execute("ng serial, nhưng chỉ 1 cổng active tại 1 thời điểm.", $true, $false, $false, $false, $false, $true, 1, $false, "ng serial, nhưng chỉ 1 cổng active tại 1 thời điểm. Thay đổi cổng bằng lệnh listen().", 2)

